$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.255.31"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "2.026.63"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'227.14"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'0.609"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'55.34"
$ws.Range("E8").Value = "  -3.97%  "
$ws.Range("D9").Value = "'0.381"
$ws.Range("E9").Value = "  -2.40%  "
$ws.Range("D10").Value = "'0.0788"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "'0.101"
$ws.Range("E11").Value = "  -5.79%  "
$ws.Range("D12").Value = "2.314.30"
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("D13").Value = "'14.32"
$ws.Range("E13").Value = "  -5.66%  "
$ws.Range("D14").Value = "'20.32"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").Value = "'0.744"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").Value = "'5.19"
$ws.Range("E16").Value = "  -3.09%  "
$ws.Range("D17").Value = "2.024.36"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").Value = "37.140.63"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "'6.33"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "'69.12"
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").Value = "0.0₃0819"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'223.59"
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'2.45"
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = "  -6.53%  "
$ws.Range("D26").Value = "'9.29"
$ws.Range("E26").Value = "  -4.92%  "
$ws.Range("D27").Value = "'165.94"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  -2.05%  "
$ws.Range("D29").Value = "'18.75"
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("D30").Value = "'1.35"
$ws.Range("E30").Value = "  -1.88%  "
$ws.Range("D31").Value = "'0.118"
$ws.Range("E31").Value = "  -2.30%  "
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "'0.0616"
$ws.Range("E33").Value = "  -1.91%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("D35").Value = "'2.36"
$ws.Range("E35").Value = "  -5.07%  "
$ws.Range("D36").Value = "'1.87"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'3.15"
$ws.Range("E38").Value = "  -4.55%  "
$ws.Range("D39").Value = "'5.54"
$ws.Range("E39").Value = "  +5.11%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.473.52"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0216"
$ws.Range("E41").Value = "  -3.40%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'16.78"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'95.93"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("D44").Value = "'0.0919"
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("E46").Value = "  -4.09%  "
$ws.Range("D47").Value = "'7.26"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("D49").Value = "'2.92"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "2.213.95"
$ws.Range("E50").Value = "  -1.76%  "
$ws.Range("D51").Value = "'3.54"
$ws.Range("E51").Value = "  -10.91%  "
